# Multiple user name and password from excel file
#
# The login test data sheet ("LoginScreen") stored a single password value
# for the "Valid Username and Valid Password Login Test" case. Update it to
# the correct/current password used for that test.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginScreen")

$ws.Range("C2").Value = "leo_123"

# The Username column now stands out as the widest column on the sheet, so
# give it its own best-fit width instead of sharing the generic column width.
$ws.Columns.Item(2).AutoFit() | Out-Null

# Leave the selection where the change was made/reviewed.
$ws.Range("C6").Select() | Out-Null
